# MENU_MOCK.xlsx edit: add a per-zone "Frequency" option (time between bleed
# ticks) for each wound zone, inserted just before that zone's StackLimit row.
#
# Zones (in sheet order): Throat, Head, Neck, Torso, Arm, Leg, Dismemberment.
# Each new row: Order=35, Type=float, Default=0.5f,
#   Tooltip="Time between bleed ticks for <zone> wounds"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row immediately above each zone's "StackLimit" row.
# Working from the bottom of the sheet upward keeps the earlier (smaller)
# row numbers valid since inserting below a row never shifts it.
$stackLimitRows = @(44, 40, 36, 32, 28, 24, 20)
foreach ($r in $stackLimitRows) {
    $ws.Rows.Item($r).Insert()
}

# After the inserts above, the new blank rows now live at:
#   20  -> Throat Frequency      (was Throat StackLimit, now pushed to 21)
#   25  -> Head Frequency        (was Head StackLimit, now pushed to 26)
#   30  -> Neck Frequency        (was Neck StackLimit, now pushed to 31)
#   35  -> Torso Frequency       (was Torso StackLimit, now pushed to 36)
#   40  -> Arm Frequency         (was Arm StackLimit, now pushed to 41)
#   45  -> Leg Frequency         (was Leg StackLimit, now pushed to 46)
#   50  -> Dismemberment Frequency (was Dismemberment StackLimit, now pushed to 51)

$newRows = @(
    @{ Row = 20; Category = "CategoryZoneThroat";        Name = "OptionThroatFrequency";        Tooltip = "Time between bleed ticks for throat wounds" },
    @{ Row = 25; Category = "CategoryZoneHead";           Name = "OptionHeadFrequency";          Tooltip = "Time between bleed ticks for head wounds" },
    @{ Row = 30; Category = "CategoryZoneNeck";           Name = "OptionNeckFrequency";          Tooltip = "Time between bleed ticks for neck wounds" },
    @{ Row = 35; Category = "CategoryZoneTorso";          Name = "OptionTorsoFrequency";         Tooltip = "Time between bleed ticks for torso wounds" },
    @{ Row = 40; Category = "CategoryZoneArm";            Name = "OptionArmFrequency";           Tooltip = "Time between bleed ticks for arm wounds" },
    @{ Row = 45; Category = "CategoryZoneLeg";            Name = "OptionLegFrequency";           Tooltip = "Time between bleed ticks for leg wounds" },
    @{ Row = 50; Category = "CategoryZoneDismemberment";  Name = "OptionDismembermentFrequency"; Tooltip = "Time between bleed ticks for dismemberment" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value = $nr.Category
    $ws.Cells.Item($r, 2).Value = 35
    $ws.Cells.Item($r, 3).Value = $nr.Name
    $ws.Cells.Item($r, 4).Value = "float"
    $ws.Cells.Item($r, 5).Value = "0.5f"
    $ws.Cells.Item($r, 6).Value = $nr.Tooltip
    # Column G ("Value Source") is left blank for these rows, matching the
    # other option rows in the sheet (it is never populated elsewhere).
}

Write-Output "done"
